$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "27.061.00"
$ws.Cells.Item(2, 5).Value = "  +0.46%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.675.38"
$ws.Cells.Item(3, 5).Value = "  +0.25%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.06%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "215.30"
$ws.Cells.Item(5, 5).Value = "  +0.25%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.17%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.08%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +1.88%  "

# Row 9
$ws.Cells.Item(9, 2).Value = "Solana"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(9, 4).Value = "21.26"
$ws.Cells.Item(9, 5).Value = "  +5.29%  "

# Row 10
$ws.Cells.Item(10, 2).Value = "Dogecoin"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(10, 4).Value = "0.0621"
$ws.Cells.Item(10, 5).Value = "  +0.27%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.75%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.912.23"
$ws.Cells.Item(12, 5).Value = "  +0.28%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.686.23"
$ws.Cells.Item(13, 5).Value = "  +0.93%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +0.88%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +1.54%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "66.00"
$ws.Cells.Item(16, 5).Value = "  +0.76%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "27.052.58"
$ws.Cells.Item(17, 5).Value = "  +0.47%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "BitcoinCash"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(18, 4).Value = "237.65"
$ws.Cells.Item(18, 5).Value = "  +1.76%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "Chainlink"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(19, 4).Value = "8.15"
$ws.Cells.Item(19, 5).Value = "  +1.30%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.0₃0742"
$ws.Cells.Item(20, 5).Value = "  +1.31%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.06%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "4.47"
$ws.Cells.Item(22, 5).Value = "  +0.96%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "9.34"
$ws.Cells.Item(23, 5).Value = "  +1.84%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "2.14"
$ws.Cells.Item(24, 5).Value = "  -1.80%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "146.61"
$ws.Cells.Item(25, 5).Value = "  +0.53%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "7.20"
$ws.Cells.Item(26, 5).Value = "  +1.14%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "16.35"
$ws.Cells.Item(27, 5).Value = "  +2.68%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.37%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -0.08%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +0.21%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Maker"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(32, 4).Value = "1.549.65"
$ws.Cells.Item(32, 5).Value = "  +6.14%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Filecoin"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(33, 4).Value = "3.35"
$ws.Cells.Item(33, 5).Value = "  +0.70%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +1.76%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "1.70"
$ws.Cells.Item(35, 5).Value = "  +3.12%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "0.598"
$ws.Cells.Item(36, 5).Value = "  +3.25%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -1.10%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +2.97%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +1.94%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +1.85%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.04%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "67.56"
$ws.Cells.Item(42, 5).Value = "  +1.95%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "5.60"
$ws.Cells.Item(43, 5).Value = "  -2.99%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -1.73%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "1.820.82"
$ws.Cells.Item(45, 5).Value = "  +0.78%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "0.782"
$ws.Cells.Item(46, 5).Value = "  +0.18%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "90.55"
$ws.Cells.Item(47, 5).Value = "  -0.06%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "0.0₆0107"
$ws.Cells.Item(48, 5).Value = "  +2.10%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "1.56"
$ws.Cells.Item(49, 5).Value = "  +1.92%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +2.54%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +5.28%  "
